# spClassification.xlsx update:
#  - Remove the "REMAINDERS" entries from column E (Sheet1). These cells
#    were the sole remaining users of that shared string, so clearing
#    them drops "REMAINDERS" from the shared-string table entirely and
#    shifts the following entry ("Coccineae") down one slot - which is
#    why D26:E32 (which held "Coccineae") end up pointing at the new,
#    lower shared-string index automatically.
#  - Add a column-width entry for column E to match its (now shorter)
#    contents.
#  - Leave the selection on H182, matching where editing finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E33:E106,E124:E130,E146:E178,E180").Areas | ForEach-Object {
    $_.ClearContents()
}

$ws.Columns.Item(5).ColumnWidth = 12.85546875

$ws.Range("H182").Select() | Out-Null
